# Daily attendance processing - 2025-12-18 20:53:47
#
# For every row in the "Recorded By" column (column G), when the list of
# recorders ends with "System" (case-insensitive), rotate the list so the
# first recorder moves to the end (e.g. "a@b.com, System" -> "System, a@b.com").
# Rows whose last recorder is something else (e.g. admin@admin.com) are left
# untouched, and single-value cells (no comma) are left untouched as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$startRow = $used.Row
$rowCount = $used.Rows.Count
$lastRow = $startRow + $rowCount - 1

$col = 7  # Column G = "Recorded By"

for ($r = $startRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "" -and $text.Contains(", ")) {
        $parts = $text -split ", "
        $last = $parts[$parts.Count - 1]

        if ($last.ToLower() -eq "system") {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = ($rotated -join ", ")
        }
    }
}
